# Updated cryptos list with latest prices/volumes scraped from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.910.67'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '2.362.93'
$ws.Range("E3").Value = '  +2.04%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '301.94'
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").Value = '95.44'
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("D7").Value = '0.505'
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '0.486'
$ws.Range("E9").Value = '  -1.48%  '
$ws.Range("D10").Value = '33.91'
$ws.Range("E10").Value = '  -1.45%  '
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("E12").Value = '  +2.86%  '
$ws.Range("D13").Value = '18.25'
$ws.Range("E13").Value = '  -4.24%  '
$ws.Range("D14").Value = '6.72'
$ws.Range("E14").Value = '  -0.28%  '
$ws.Range("D15").Value = '2.730.32'
$ws.Range("E15").Value = '  +1.92%  '
$ws.Range("D16").Value = '2.345.36'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").Value = '0.796'
$ws.Range("E17").Value = '  +0.72%  '
$ws.Range("D18").Value = '42.855.85'
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("D19").Value = '12.04'
$ws.Range("E19").Value = '  -0.86%  '
$ws.Range("E20").Value = '  +1.79%  '
$ws.Range("D21").Value = '0.0₃0886'
$ws.Range("E21").Value = '  -0.63%  '
$ws.Range("D22").Value = '67.84'
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").Value = '234.99'
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("E24").Value = '  -2.52%  '
$ws.Range("E26").Value = '  +0.22%  '
$ws.Range("D27").Value = '24.66'
$ws.Range("E27").Value = '  +1.00%  '
$ws.Range("E28").Value = '  +0.29%  '
$ws.Range("D29").Value = '9.24'
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("D30").Value = '31.51'
$ws.Range("E30").Value = '  -2.13%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").Value = '5.04'
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("B33").Value = 'Celestia'
$ws.Range("C33").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D33").Value = '17.32'
$ws.Range("E33").Value = '  -3.05%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.0724'
$ws.Range("E34").Value = '  +3.67%  '
$ws.Range("D35").Value = '0.105'
$ws.Range("E35").Value = '  +4.70%  '
$ws.Range("E36").Value = '  +3.53%  '
$ws.Range("D37").Value = '4.36'
$ws.Range("E38").Value = '  -0.98%  '
$ws.Range("D39").Value = '2.79'
$ws.Range("E39").Value = '  +2.28%  '
$ws.Range("D40").Value = '120.71'
$ws.Range("E40").Value = '  -27.52%  '
$ws.Range("E41").Value = '  -0.55%  '
$ws.Range("D42").Value = '21.40'
$ws.Range("E42").Value = '  +2.87%  '
$ws.Range("D43").Value = '1.933.38'
$ws.Range("E43").Value = '  +0.27%  '
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("E45").Value = '  +2.07%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '9.19'
$ws.Range("E46").Value = '  -9.28%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '2.71'
$ws.Range("E47").Value = '  -1.83%  '
$ws.Range("D48").Value = '2.591.61'
$ws.Range("E48").Value = '  +1.79%  '
$ws.Range("E49").Value = '  +1.95%  '
$ws.Range("D50").Value = '72.02'
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("E51").Value = '  +1.16%  '
